# Applies the cryptos-list refresh described by the commit:
# "Updated cryptos list on Thu Aug 15 18:56:24 UTC 2024 with GitHub Actions"
#
# Source data cells are plain text (t="inlineStr") even when the text looks
# like a number (e.g. "517.66"), so every write below keeps the destination
# cell as Text. For values that parse as a number, Excel's COM layer would
# otherwise silently coerce the write to a numeric cell, so those go through
# NumberFormat "@" first and have their style reset back to Normal afterwards
# so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, price/volume strings that are not
# ambiguous with a number).
$textUpdates = @(
    @{ Cell = "D2"; Value = '57.760.78' },
    @{ Cell = "E2"; Value = '  -2.00%  ' },
    @{ Cell = "D3"; Value = '2.544.11' },
    @{ Cell = "E3"; Value = '  -4.04%  ' },
    @{ Cell = "E4"; Value = '  +0.08%  ' },
    @{ Cell = "E5"; Value = '  -1.14%  ' },
    @{ Cell = "E6"; Value = '  -4.22%  ' },
    @{ Cell = "E7"; Value = '  +0.12%  ' },
    @{ Cell = "E8"; Value = '  -1.42%  ' },
    @{ Cell = "E9"; Value = '  -7.08%  ' },
    @{ Cell = "E10"; Value = '  -3.71%  ' },
    @{ Cell = "E11"; Value = '  -3.01%  ' },
    @{ Cell = "E12"; Value = '  -0.03%  ' },
    @{ Cell = "D13"; Value = '2.997.53' },
    @{ Cell = "E13"; Value = '  -3.71%  ' },
    @{ Cell = "D14"; Value = '57.747.79' },
    @{ Cell = "E14"; Value = '  -2.07%  ' },
    @{ Cell = "E15"; Value = '  -5.09%  ' },
    @{ Cell = "E16"; Value = '  -3.03%  ' },
    @{ Cell = "D17"; Value = '2.538.95' },
    @{ Cell = "E17"; Value = '  -3.70%  ' },
    @{ Cell = "E18"; Value = '  -2.04%  ' },
    @{ Cell = "E19"; Value = '  -2.14%  ' },
    @{ Cell = "E20"; Value = '  -2.47%  ' },
    @{ Cell = "E21"; Value = '  -3.85%  ' },
    @{ Cell = "E22"; Value = '  -0.08%  ' },
    @{ Cell = "E23"; Value = '  +1.61%  ' },
    @{ Cell = "E24"; Value = '  -0.64%  ' },
    @{ Cell = "B25"; Value = 'Binance-PegBSC-USD' },
    @{ Cell = "C25"; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd' },
    @{ Cell = "E25"; Value = '  -0.04%  ' },
    @{ Cell = "B26"; Value = 'Polygon' },
    @{ Cell = "C26"; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Cell = "E26"; Value = '  -4.32%  ' },
    @{ Cell = "D27"; Value = '2.678.38' },
    @{ Cell = "E27"; Value = '  -3.37%  ' },
    @{ Cell = "E28"; Value = '  -2.55%  ' },
    @{ Cell = "D29"; Value = '0.0₃0750' },
    @{ Cell = "E29"; Value = '  -6.36%  ' },
    @{ Cell = "E30"; Value = '  +0.08%  ' },
    @{ Cell = "E31"; Value = '  -7.44%  ' },
    @{ Cell = "E32"; Value = '  -1.89%  ' },
    @{ Cell = "E33"; Value = '  -0.44%  ' },
    @{ Cell = "E34"; Value = '  -1.98%  ' },
    @{ Cell = "E35"; Value = '  -4.66%  ' },
    @{ Cell = "E36"; Value = '  -5.64%  ' },
    @{ Cell = "E37"; Value = '  -6.61%  ' },
    @{ Cell = "E38"; Value = '  -2.93%  ' },
    @{ Cell = "E39"; Value = '  -5.69%  ' },
    @{ Cell = "E40"; Value = '  -5.00%  ' },
    @{ Cell = "E41"; Value = '  +0.13%  ' },
    @{ Cell = "E42"; Value = '  -3.52%  ' },
    @{ Cell = "E43"; Value = '  -0.08%  ' },
    @{ Cell = "E44"; Value = '  -1.54%  ' },
    @{ Cell = "E45"; Value = '  -6.52%  ' },
    @{ Cell = "E47"; Value = '  -2.66%  ' },
    @{ Cell = "E48"; Value = '  -7.02%  ' },
    @{ Cell = "D49"; Value = '1.976.39' },
    @{ Cell = "E49"; Value = '  -2.74%  ' },
    @{ Cell = "E50"; Value = '  -3.31%  ' },
    @{ Cell = "E51"; Value = '  -5.46%  ' }
)
foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Price updates whose text would otherwise be auto-detected as a number by
# the COM layer (e.g. "0.400", "65.00", "260.00") -- force Text explicitly so
# the stored value/type matches the source (and trailing zeros survive).
$numericLookingUpdates = @(
    @{ Cell = "D5"; Value = '517.66' },
    @{ Cell = "D6"; Value = '138.42' },
    @{ Cell = "D7"; Value = '1.00' },
    @{ Cell = "D9"; Value = '6.51' },
    @{ Cell = "D10"; Value = '0.0989' },
    @{ Cell = "D15"; Value = '19.96' },
    @{ Cell = "D18"; Value = '332.86' },
    @{ Cell = "D19"; Value = '4.27' },
    @{ Cell = "D20"; Value = '10.09' },
    @{ Cell = "D21"; Value = '6.11' },
    @{ Cell = "D23"; Value = '65.00' },
    @{ Cell = "D25"; Value = '0.999' },
    @{ Cell = "D26"; Value = '0.400' },
    @{ Cell = "D28"; Value = '6.90' },
    @{ Cell = "D31"; Value = '6.18' },
    @{ Cell = "D32"; Value = '1.56' },
    @{ Cell = "D33"; Value = '148.89' },
    @{ Cell = "D34"; Value = '18.44' },
    @{ Cell = "D35"; Value = '3.95' },
    @{ Cell = "D36"; Value = '1.12' },
    @{ Cell = "D38"; Value = '35.69' },
    @{ Cell = "D44"; Value = '0.0952' },
    @{ Cell = "D45"; Value = '0.575' },
    @{ Cell = "D46"; Value = '260.00' },
    @{ Cell = "D47"; Value = '0.0518' },
    @{ Cell = "D48"; Value = '18.52' },
    @{ Cell = "D51"; Value = '4.50' }
)
foreach ($u in $numericLookingUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
